$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The node ids in column A for the existing rows 20-39 shift down by 8 to
# close the numbering gap, making room for the two new "escala" nodes that
# get appended at the end of the table.
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21
$ws.Range("A24").Value = 22
$ws.Range("A25").Value = 23
$ws.Range("A26").Value = 24
$ws.Range("A27").Value = 25
$ws.Range("A28").Value = 26
$ws.Range("A29").Value = 27
$ws.Range("A30").Value = 28
$ws.Range("A31").Value = 29
$ws.Range("A32").Value = 30
$ws.Range("A33").Value = 31
$ws.Range("A34").Value = 32
$ws.Range("A35").Value = 33
$ws.Range("A36").Value = 34
$ws.Range("A37").Value = 35
$ws.Range("A38").Value = 36
$ws.Range("A39").Value = 37

# Agrega los dos nodos de escala junto con su enlace.
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "escala1"
$ws.Range("C40").Value = 9
$ws.Range("D40").Value = 1

$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "escala2"
$ws.Range("C41").Value = 10
$ws.Range("D41").Value = 1

# Leave the cursor on the last edited link, matching the author's final
# selection after adding the new rows.
$ws.Range("C37").Select()

# Columns picked up a hair-narrower "optimal width" once the new rows were
# added; nudge them to the closest width this engine can represent.
$ws.Columns.Item(1).ColumnWidth = 7.43
$ws.Columns.Item(2).ColumnWidth = 16.75
